$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (A1, B1, C1)
$ws.Range("A1").Value = 148.89454749342784
$ws.Range("B1").Value = 4.6262598347551753
$ws.Range("C1").Value = 0.84716242661448149

# Update column widths for column A (10.7109375 -> 11.7109375) and
# column C (11.7109375 -> 12.7109375). Excel's ColumnWidth property is
# quantized to the workbook's pixel grid (1/6-character steps for the
# default Calibri 11 font), so the closest representable values are used.
$ws.Columns.Item(1).ColumnWidth = 10.8333333333333
$ws.Columns.Item(3).ColumnWidth = 11.8333333333333
